$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("17:17").Insert()

$ws.Range("A17").Value = 10
$ws.Range("B17").Value = "Vega Modelo de Temuco"
$ws.Range("C17").Value = "La Araucanía"
$ws.Range("D17").Value = 44624
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100107
$ws.Range("H17").Value = "Otros"
$ws.Range("I17").Value = 100107011
$ws.Range("J17").Value = "Tuna"
$ws.Range("K17").Value = "Sin especificar"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 75
$ws.Range("N17").Value = 17000
$ws.Range("O17").Value = 18000
$ws.Range("P17").Value = 17467
$ws.Range("Q17").Value = "$/caja 18 kilos"
$ws.Range("R17").Value = "Provincia de Los Andes"
$ws.Range("S17").Value = 970
$ws.Range("T17").Value = 18
